# Apply the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain decimal number (e.g. "581.52") need their
# NumberFormat forced to text ("@") first, otherwise Excel auto-converts the
# literal to a numeric value instead of keeping it as the original inline string.
$textCells = @('D5', 'D6', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D34', 'D36', 'D37', 'D40', 'D41', 'D44', 'D45', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.920.23'
$ws.Range('E2').Value = '  +5.75%  '
$ws.Range('D3').Value = '2.978.44'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('D5').Value = '581.52'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').Value = '153.42'
$ws.Range('E6').Value = '  +7.13%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = '2.976.24'
$ws.Range('E9').Value = '  +3.03%  '
$ws.Range('E10').Value = '  +3.25%  '
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '34.00'
$ws.Range('E14').Value = '  +6.46%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '64.849.20'
$ws.Range('E16').Value = '  +5.60%  '
$ws.Range('D17').Value = '3.468.45'
$ws.Range('E17').Value = '  +2.91%  '
$ws.Range('D18').Value = '6.91'
$ws.Range('E18').Value = '  +4.53%  '
$ws.Range('D19').Value = '2.972.67'
$ws.Range('E19').Value = '  +2.92%  '
$ws.Range('D20').Value = '448.65'
$ws.Range('E20').Value = '  +3.38%  '
$ws.Range('D21').Value = '13.68'
$ws.Range('E21').Value = '  +3.94%  '
$ws.Range('D22').Value = '0.680'
$ws.Range('E22').Value = '  +3.72%  '
$ws.Range('D23').Value = '7.26'
$ws.Range('E23').Value = '  +5.94%  '
$ws.Range('D24').Value = '80.96'
$ws.Range('E24').Value = '  +2.21%  '
$ws.Range('D25').Value = '12.29'
$ws.Range('E25').Value = '  +3.95%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').Value = '  +8.84%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '10.59'
$ws.Range('E27').Value = '  +5.84%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').Value = '7.83'
$ws.Range('E29').Value = '  +11.50%  '
$ws.Range('D30').Value = '2.37'
$ws.Range('E30').Value = '  +15.04%  '
$ws.Range('D31').Value = '2.58'
$ws.Range('E31').Value = '  +3.29%  '
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('E33').Value = '  +3.73%  '
$ws.Range('D34').Value = '26.74'
$ws.Range('E34').Value = '  +4.43%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '0.979'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').Value = '5.68'
$ws.Range('E37').Value = '  +4.85%  '
$ws.Range('E38').Value = '  +7.64%  '
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('B40').Value = 'Arweave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D40').Value = '44.34'
$ws.Range('E40').Value = '  +12.09%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('E42').Value = '  +5.36%  '
$ws.Range('E43').Value = '  +10.57%  '
$ws.Range('D44').Value = '8.40'
$ws.Range('E44').Value = '  +1.84%  '
$ws.Range('D45').Value = '381.96'
$ws.Range('E45').Value = '  +12.12%  '
$ws.Range('D46').Value = '2.763.20'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('D47').Value = '0.0348'
$ws.Range('E47').Value = '  +4.19%  '
$ws.Range('D48').Value = '134.38'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '23.12'
$ws.Range('E50').Value = '  +7.43%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.105'
$ws.Range('E51').Value = '  +2.15%  '
